$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3: Kinect purchase amount correction, and clear old Display Board / 67 entries (F3/G3) ----
$ws.Range("D3").Value = 146.89
$ws.Range("F3").ClearContents()
$ws.Range("G3").ClearContents()

# ---- Row 5: was "Printing (Rowan)" -> now "PCBs"; amount updated; note cleared ----
$ws.Range("A5").Value = "PCBs"
$ws.Range("D5").Value = 144.07
$ws.Range("E5").ClearContents()

# ---- Row 6: Sheet Metal / Rowan H (labels unchanged); clear the old funding entry (J6/K6) ----
$ws.Range("A6").Value = "Sheet Metal"
$ws.Range("B6").Value = "Rowan H"
$ws.Range("J6").ClearContents()
$ws.Range("K6").ClearContents()

# ---- Row 7: Kinect Mount / Rowan H (labels unchanged) ----
$ws.Range("A7").Value = "Kinect Mount"
$ws.Range("B7").Value = "Rowan H"

# ---- Row 8: Sonar Customs / Iain P (unchanged) ----
$ws.Range("A8").Value = "Sonar Customs"
$ws.Range("B8").Value = "Iain P"

# ---- Row 9: PCB Customs, now purchased by Jordan V, new amount ----
$ws.Range("A9").Value = "PCB Customs"
$ws.Range("B9").Value = "Jordan V"
$ws.Range("D9").Value = 21.65

# ---- Row 10: Paint, now purchased by Jordan V, new amount ----
$ws.Range("A10").Value = "Paint"
$ws.Range("B10").Value = "Jordan V"
$ws.Range("D10").Value = 17.05

# ---- Row 11: Parts for PCB, now purchased by Jordan V, amount unchanged ----
$ws.Range("A11").Value = "Parts for PCB"
$ws.Range("B11").Value = "Jordan V"
$ws.Range("D11").Value = 85.76

# ---- Row 12 (new): Display Board, purchased by Rowan H, moved from F3/G3 ----
$ws.Range("A12").Value = "Display Board"
$ws.Range("B12").Value = "Rowan H"
$ws.Range("D12").Value = 67

# ---- New shared strings must be introduced in this exact order to match the
#      canonical shared-string table layout: Jordan, Rowan, Iain, Per Person, Refund ----
$ws.Range("A20").Value = "Jordan"
$ws.Range("A21").Value = "Rowan"
$ws.Range("A22").Value = "Iain"
$ws.Range("E19").Value = "Per Person"
$ws.Range("C19").Value = "Refund"

# ---- Row 19 (new): Refund / Per Person formatting + per-person refund formula ----
$ws.Range("C19").Font.Bold = $true
$ws.Range("E19").Font.Bold = $true
$ws.Range("F19").Formula = "=F18/3"

# ---- Row 20 (new): Jordan's totals ----
$ws.Range("B20").Formula = "=D3+D5+D9+D10+D11"
$ws.Range("C20").Formula = "=B20-F19"

# ---- Row 21 (new): Rowan's totals ----
$ws.Range("B21").Formula = "=D6+D7+D12"
$ws.Range("C21").Formula = "=B21-F19"

# ---- Row 22 (new): Iain's totals ----
$ws.Range("B22").Formula = "=D4+D8"
$ws.Range("C22").Formula = "=B22-F19"

# ---- Update the active selection to mirror the committed workbook state ----
$ws.Range("D21").Select()
